$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.899.24"
$ws.Range("E2").Value = "  -2.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.578.59"
$ws.Range("E3").Value = "  -4.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.09"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.49"
$ws.Range("E6").Value = "  -2.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +1.16%  "

# Row 9
$ws.Range("E9").Value = "  -1.70%  "

# Row 10
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("E11").Value = "  +3.20%  "

# Row 12
$ws.Range("E12").Value = "  -0.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.034.74"
$ws.Range("E13").Value = "  -3.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.42"
$ws.Range("E14").Value = "  -3.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.844.60"
$ws.Range("E15").Value = "  -2.29%  "

# Row 16
$ws.Range("E16").Value = "  -1.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.577.21"
$ws.Range("E17").Value = "  -4.03%  "

# Row 18
$ws.Range("E18").Value = "  -3.74%  "

# Row 19
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.18"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21
$ws.Range("E21").Value = "  -4.99%  "

# Row 22
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.490"
$ws.Range("E23").Value = "  -3.71%  "

# Row 24
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$ws.Range("E25").Value = "  -1.26%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.12"
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0832"
$ws.Range("E28").Value = "  -3.29%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  +2.45%  "

# Row 30
$ws.Range("E30").Value = "  -1.28%  "

# Row 31
$ws.Range("E31").Value = "  -3.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.63"
$ws.Range("E32").Value = "  -1.95%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.86"
$ws.Range("E33").Value = "  +0.63%  "

# Row 34
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.12"
$ws.Range("E35").Value = "  -2.40%  "

# Row 36
$ws.Range("E36").Value = "  -2.69%  "

# Row 37
$ws.Range("E37").Value = "  -0.65%  "

# Row 38
$ws.Range("E38").Value = "  -2.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "325.77"
$ws.Range("E39").Value = "  -5.19%  "

# Row 40
$ws.Range("E40").Value = "  -4.77%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  -0.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.50"
$ws.Range("E42").Value = "  -1.65%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.58"
$ws.Range("E43").Value = "  -1.34%  "

# Row 44
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").Value = "  -2.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.93"
$ws.Range("E46").Value = "  -1.11%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.55"
$ws.Range("E47").Value = "  -4.05%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0545"
$ws.Range("E48").Value = "  -3.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0965"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50
$ws.Range("E50").Value = "  -1.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.047.60"
$ws.Range("E51").Value = "  -2.10%  "
